$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" contain the same data table and both need
# their "想去人数" (F column) counts bumped for rows 2 and 3.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 688
    $ws.Range("F3").Value = 4003
}
